$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.378.96"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.885.98"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'237.98"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").Value = "'0.4688"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").Value = "'0.2820"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").Value = "'0.06578"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("D10").Value = "'19.74"
$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("D11").Value = "'98.55"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").Value = "'0.07737"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").Value = "1.880.67"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "'5.143"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "'0.6693"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "'285.17"
$ws.Range("E16").Value = "  +12.37%  "

$ws.Range("D17").Value = "30.345.44"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "'0.9983"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "'12.63"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").Value = "2.131.09"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").Value = "'0.000007311"
$ws.Range("E21").Value = "  -1.65%  "

$ws.Range("D22").Value = "'5.339"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "'0.9985"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("B24").Value = "BitDAO"
$ws.Range("C24").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D24").Value = "'0.4584"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'6.200"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "'9.299"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'167.08"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.07"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.990"
$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.376"
$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.09848"
$ws.Range("E31").Value = "  -2.18%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.473"
$ws.Range("E32").Value = "  -4.22%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.496"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.194"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.04703"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7095"
$ws.Range("E36").Value = "  -1.83%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.099"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'0.9982"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.703"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01873"
$ws.Range("E40").Value = "  -1.24%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.679"
$ws.Range("E41").Value = "  +7.87%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.526"
$ws.Range("E42").Value = "  -2.42%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'72.77"
$ws.Range("E43").Value = "  +0.54%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8699"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.967"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'104.24"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'0.9984"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4210"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'993.82"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.447"
$ws.Range("E50").Value = "  +8.80%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.260"
$ws.Range("E51").Value = "  -1.39%  "
